{"js": "const REPLACEMENTS = [\n  [\"2023-05-16 Tuesday\", \"2023-05-17 Wednesday\"],\n  [\"78\u00d761=4758\", \"79\u00d727=2133\"],\n  [\"33\u00d754=1782\", \"90\u00d748=4320\"],\n  [\"29\u00d717=493\", \"66\u00d717=1122\"],\n  [\"53\u00d717=901\", \"96\u00d737=3552\"],\n  [\"57\u00d748=2736\", \"82\u00d748=3936\"],\n  [\"25\u00d774=1850\", \"83\u00d713=1079\"],\n  [\"98\u00d760=5880\", \"12\u00d781=972\"],\n  [\"94\u00d719=1786\", \"91\u00d791=8281\"],\n  [\"44\u00d795=4180\", \"34\u00d752=1768\"],\n  [\"45\u00d741=1845\", \"63\u00d785=5355\"],\n  [\"40\u00d712=480\", \"77\u00d768=5236\"],\n  [\"31\u00d712=372\", \"94\u00d775=7050\"],\n  [\"14\u00d724=336\", \"18\u00d729=522\"],\n  [\"72\u00d746=3312\", \"79\u00d785=6715\"],\n  [\"22\u00d775=1650\", \"86\u00d737=3182\"],\n  [\"17\u00d737=629\", \"40\u00d742=1680\"],\n  [\"65\u00d759=3835\", \"26\u00d729=754\"],\n  [\"23\u00d738=874\", \"17\u00d780=1360\"],\n  [\"47\u00d729=1363\", \"77\u00d737=2849\"],\n  [\"67\u00d745=3015\", \"92\u00d761=5612\"],\n  [\"84\u00d757=4788\", \"86\u00d763=5418\"],\n  [\"90\u00d781=7290\", \"16\u00d796=1536\"],\n  [\"66\u00d779=5214\", \"70\u00d719=1330\"],\n  [\"18\u00d781=1458\", \"19\u00d753=1007\"],\n  [\"84\u00d729=2436\", \"17\u00d746=782\"],\n  [\"61\u00d785=5185\", \"100\u00d765=6500\"],\n  [\"57\u00d757=3249\", \"78\u00d795=7410\"],\n  [\"33\u00d710=330\", \"30\u00d781=2430\"],\n  [\"33\u00d781=2673\", \"43\u00d777=3311\"],\n  [\"32\u00d773=2336\", \"67\u00d738=2546\"],\n  [\"96\u00d762=5952\", \"33\u00d723=759\"],\n  [\"100\u00d776=7600\", \"95\u00d723=2185\"],\n  [\"21\u00d729=609\", \"39\u00d755=2145\"],\n  [\"93\u00d795=8835\", \"44\u00d724=1056\"],\n  [\"70\u00d776=5320\", \"54\u00d791=4914\"],\n  [\"59\u00d741=2419\", \"23\u00d712=276\"],\n  [\"31\u00d773=2263\", \"60\u00d720=1200\"],\n  [\"23\u00d750=1150\", \"25\u00d791=2275\"],\n  [\"47\u00d752=2444\", \"23\u00d755=1265\"],\n  [\"38\u00d795=3610\", \"32\u00d768=2176\"],\n  [\"42\u00d729=1218\", \"36\u00d710=360\"],\n  [\"16\u00d770=1120\", \"43\u00d798=4214\"],\n  [\"83\u00d715=1245\", \"96\u00d717=1632\"],\n  [\"60\u00d724=1440\", \"75\u00d743=3225\"],\n  [\"66\u00d758=3828\", \"71\u00d749=3479\"],\n  [\"40\u00d739=1560\", \"30\u00d791=2730\"],\n  [\"20\u00d752=1040\", \"67\u00d719=1273\"],\n  [\"38\u00d751=1938\", \"19\u00d793=1767\"],\n  [\"20\u00d732=640\", \"94\u00d722=2068\"],\n  [\"35\u00d720=700\", \"32\u00d738=1216\"],\n  [\"37\u00d743=1591\", \"96\u00d715=1440\"],\n  [\"42\u00d720=840\", \"37\u00d740=1480\"],\n  [\"92\u00d727=2484\", \"14\u00d739=546\"],\n  [\"58\u00d734=1972\", \"100\u00d741=4100\"],\n  [\"93\u00d747=4371\", \"76\u00d743=3268\"],\n  [\"80\u00d774=5920\", \"58\u00d798=5684\"],\n  [\"91\u00d763=5733\", \"27\u00d787=2349\"],\n  [\"70\u00d718=1260\", \"23\u00d711=253\"],\n  [\"96\u00d778=7488\", \"86\u00d734=2924\"],\n  [\"37\u00d782=3034\", \"91\u00d786=7826\"],\n  [\"82\u00d784=6888\", \"32\u00d780=2560\"],\n  [\"76\u00d797=7372\", \"96\u00d714=1344\"],\n  [\"96\u00d727=2592\", \"66\u00d790=5940\"],\n  [\"61\u00d737=2257\", \"72\u00d735=2520\"],\n  [\"76\u00d730=2280\", \"89\u00d758=5162\"],\n  [\"99\u00d716=1584\", \"100\u00d740=4000\"],\n  [\"87\u00d783=7221\", \"35\u00d785=2975\"],\n  [\"93\u00d754=5022\", \"19\u00d731=589\"],\n  [\"55\u00d787=4785\", \"36\u00d722=792\"],\n  [\"12\u00d758=696\", \"63\u00d754=3402\"],\n  [\"30\u00d795=2850\", \"25\u00d762=1550\"],\n  [\"27\u00d741=1107\", \"81\u00d764=5184\"],\n  [\"54\u00d752=2808\", \"58\u00d738=2204\"],\n  [\"91\u00d769=6279\", \"48\u00d758=2784\"],\n  [\"60\u00d744=2640\", \"65\u00d772=4680\"],\n  [\"88\u00d785=7480\", \"58\u00d785=4930\"],\n  [\"28\u00d752=1456\", \"62\u00d757=3534\"],\n  [\"100\u00d790=9000\", \"26\u00d798=2548\"],\n  [\"95\u00d788=8360\", \"93\u00d742=3906\"],\n  [\"25\u00d757=1425\", \"56\u00d750=2800\"],\n  [\"23\u00d772=1656\", \"83\u00d765=5395\"],\n  [\"36\u00d781=2916\", \"12\u00d747=564\"],\n  [\"28\u00d712=336\", \"29\u00d781=2349\"],\n  [\"43\u00d764=2752\", \"22\u00d750=1100\"],\n  [\"43\u00d758=2494\", \"14\u00d766=924\"],\n  [\"54\u00d710=540\", \"53\u00d724=1272\"],\n  [\"17\u00d775=1275\", \"90\u00d740=3600\"],\n  [\"60\u00d785=5100\", \"73\u00d765=4745\"],\n  [\"63\u00d771=4473\", \"21\u00d774=1554\"],\n  [\"53\u00d725=1325\", \"58\u00d779=4582\"],\n  [\"92\u00d789=8188\", \"91\u00d773=6643\"],\n  [\"34\u00d724=816\", \"52\u00d771=3692\"],\n  [\"56\u00d798=5488\", \"17\u00d774=1258\"],\n  [\"23\u00d788=2024\", \"51\u00d722=1122\"],\n  [\"97\u00d753=5141\", \"54\u00d753=2862\"],\n  [\"11\u00d715=165\", \"30\u00d760=1800\"],\n  [\"43\u00d716=688\", \"44\u00d734=1496\"],\n  [\"89\u00d754=4806\", \"24\u00d732=768\"],\n  [\"96\u00d749=4704\", \"63\u00d760=3780\"],\n  [\"29\u00d721=609\", \"12\u00d783=996\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of REPLACEMENTS) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$Replacements = @(\n    @('2023-05-16 Tuesday', '2023-05-17 Wednesday'),\n    @('78\u00d761=4758', '79\u00d727=2133'),\n    @('33\u00d754=1782', '90\u00d748=4320'),\n    @('29\u00d717=493', '66\u00d717=1122'),\n    @('53\u00d717=901', '96\u00d737=3552'),\n    @('57\u00d748=2736', '82\u00d748=3936'),\n    @('25\u00d774=1850', '83\u00d713=1079'),\n    @('98\u00d760=5880', '12\u00d781=972'),\n    @('94\u00d719=1786', '91\u00d791=8281'),\n    @('44\u00d795=4180', '34\u00d752=1768'),\n    @('45\u00d741=1845', '63\u00d785=5355'),\n    @('40\u00d712=480', '77\u00d768=5236'),\n    @('31\u00d712=372', '94\u00d775=7050'),\n    @('14\u00d724=336', '18\u00d729=522'),\n    @('72\u00d746=3312', '79\u00d785=6715'),\n    @('22\u00d775=1650', '86\u00d737=3182'),\n    @('17\u00d737=629', '40\u00d742=1680'),\n    @('65\u00d759=3835', '26\u00d729=754'),\n    @('23\u00d738=874', '17\u00d780=1360'),\n    @('47\u00d729=1363', '77\u00d737=2849'),\n    @('67\u00d745=3015', '92\u00d761=5612'),\n    @('84\u00d757=4788', '86\u00d763=5418'),\n    @('90\u00d781=7290', '16\u00d796=1536'),\n    @('66\u00d779=5214', '70\u00d719=1330'),\n    @('18\u00d781=1458', '19\u00d753=1007'),\n    @('84\u00d729=2436', '17\u00d746=782'),\n    @('61\u00d785=5185', '100\u00d765=6500'),\n    @('57\u00d757=3249', '78\u00d795=7410'),\n    @('33\u00d710=330', '30\u00d781=2430'),\n    @('33\u00d781=2673', '43\u00d777=3311'),\n    @('32\u00d773=2336', '67\u00d738=2546'),\n    @('96\u00d762=5952', '33\u00d723=759'),\n    @('100\u00d776=7600', '95\u00d723=2185'),\n    @('21\u00d729=609', '39\u00d755=2145'),\n    @('93\u00d795=8835', '44\u00d724=1056'),\n    @('70\u00d776=5320', '54\u00d791=4914'),\n    @('59\u00d741=2419', '23\u00d712=276'),\n    @('31\u00d773=2263', '60\u00d720=1200'),\n    @('23\u00d750=1150', '25\u00d791=2275'),\n    @('47\u00d752=2444', '23\u00d755=1265'),\n    @('38\u00d795=3610', '32\u00d768=2176'),\n    @('42\u00d729=1218', '36\u00d710=360'),\n    @('16\u00d770=1120', '43\u00d798=4214'),\n    @('83\u00d715=1245', '96\u00d717=1632'),\n    @('60\u00d724=1440', '75\u00d743=3225'),\n    @('66\u00d758=3828', '71\u00d749=3479'),\n    @('40\u00d739=1560', '30\u00d791=2730'),\n    @('20\u00d752=1040', '67\u00d719=1273'),\n    @('38\u00d751=1938', '19\u00d793=1767'),\n    @('20\u00d732=640', '94\u00d722=2068'),\n    @('35\u00d720=700', '32\u00d738=1216'),\n    @('37\u00d743=1591', '96\u00d715=1440'),\n    @('42\u00d720=840', '37\u00d740=1480'),\n    @('92\u00d727=2484', '14\u00d739=546'),\n    @('58\u00d734=1972', '100\u00d741=4100'),\n    @('93\u00d747=4371', '76\u00d743=3268'),\n    @('80\u00d774=5920', '58\u00d798=5684'),\n    @('91\u00d763=5733', '27\u00d787=2349'),\n    @('70\u00d718=1260', '23\u00d711=253'),\n    @('96\u00d778=7488', '86\u00d734=2924'),\n    @('37\u00d782=3034', '91\u00d786=7826'),\n    @('82\u00d784=6888', '32\u00d780=2560'),\n    @('76\u00d797=7372', '96\u00d714=1344'),\n    @('96\u00d727=2592', '66\u00d790=5940'),\n    @('61\u00d737=2257', '72\u00d735=2520'),\n    @('76\u00d730=2280', '89\u00d758=5162'),\n    @('99\u00d716=1584', '100\u00d740=4000'),\n    @('87\u00d783=7221', '35\u00d785=2975'),\n    @('93\u00d754=5022', '19\u00d731=589'),\n    @('55\u00d787=4785', '36\u00d722=792'),\n    @('12\u00d758=696', '63\u00d754=3402'),\n    @('30\u00d795=2850', '25\u00d762=1550'),\n    @('27\u00d741=1107', '81\u00d764=5184'),\n    @('54\u00d752=2808', '58\u00d738=2204'),\n    @('91\u00d769=6279', '48\u00d758=2784'),\n    @('60\u00d744=2640', '65\u00d772=4680'),\n    @('88\u00d785=7480', '58\u00d785=4930'),\n    @('28\u00d752=1456', '62\u00d757=3534'),\n    @('100\u00d790=9000', '26\u00d798=2548'),\n    @('95\u00d788=8360', '93\u00d742=3906'),\n    @('25\u00d757=1425', '56\u00d750=2800'),\n    @('23\u00d772=1656', '83\u00d765=5395'),\n    @('36\u00d781=2916', '12\u00d747=564'),\n    @('28\u00d712=336', '29\u00d781=2349'),\n    @('43\u00d764=2752', '22\u00d750=1100'),\n    @('43\u00d758=2494', '14\u00d766=924'),\n    @('54\u00d710=540', '53\u00d724=1272'),\n    @('17\u00d775=1275', '90\u00d740=3600'),\n    @('60\u00d785=5100', '73\u00d765=4745'),\n    @('63\u00d771=4473', '21\u00d774=1554'),\n    @('53\u00d725=1325', '58\u00d779=4582'),\n    @('92\u00d789=8188', '91\u00d773=6643'),\n    @('34\u00d724=816', '52\u00d771=3692'),\n    @('56\u00d798=5488', '17\u00d774=1258'),\n    @('23\u00d788=2024', '51\u00d722=1122'),\n    @('97\u00d753=5141', '54\u00d753=2862'),\n    @('11\u00d715=165', '30\u00d760=1800'),\n    @('43\u00d716=688', '44\u00d734=1496'),\n    @('89\u00d754=4806', '24\u00d732=768'),\n    @('96\u00d749=4704', '63\u00d760=3780'),\n    @('29\u00d721=609', '12\u00d783=996'),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $Replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $result = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $result) {\n        throw \"Find/Replace failed for: $oldText\"\n    }\n}\n\nWrite-Output \"done\"\n"}
